# Update "want to go" (想去人数) counts in both the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 1327
    $ws.Range("F3").Value = 1826
    $ws.Range("F4").Value = 138
    $ws.Range("F6").Value = 6294
    $ws.Range("F7").Value = 155
}
